$d = $word.ActiveDocument
$sel = $word.Selection
$sel.EndKey(6)
$sel.TypeText("`r`rWorking for lms project.")
Write-Host "typed"
foreach ($p in $d.Paragraphs) {
    Write-Host "Para: [$($p.Range.Text)]"
}
# Now try replacing "singh." with itself to see if bookmark relocates due to range edit tracking
$d.Content.Find.Execute("singh.", $true, $false, $false, $false, $false, $true, 1, $false, "singh!", 2)
Write-Host $d.Content.WordOpenXML
